$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.482.44"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.54"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6934"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3063"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.60"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07771"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.27"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.138"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6904"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.74"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.401"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.457.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008264"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.100.77"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.08"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.664"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1492"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.894"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.01%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.73"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.26"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.536"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.247"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.157"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.195"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05106"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7702"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.890"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.150"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.684"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.335.14"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01866"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9732"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.59%  "

$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.33"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.849"
$ws.Range("D43").ClearFormats()

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.766"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.999.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5219"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.782"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.17"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.973"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.54%  "
